# Apply updated "想去人数" (want-to-go count) values across the four
# worksheets of the workbook, matching the data refresh captured in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 871
$ws1.Range("F3").Value = 13859
$ws1.Range("F4").Value = 13647
$ws1.Range("F5").Value = 1055
$ws1.Range("F8").Value = 605
$ws1.Range("F12").Value = 772
$ws1.Range("F13").Value = 2154
$ws1.Range("F14").Value = 117
$ws1.Range("F15").Value = 94
$ws1.Range("F16").Value = 80
$ws1.Range("F17").Value = 130
$ws1.Range("F19").Value = 538
$ws1.Range("F20").Value = 438
$ws1.Range("F21").Value = 417
$ws1.Range("F22").Value = 331
$ws1.Range("F24").Value = 845
$ws1.Range("F25").Value = 100

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 169
$ws2.Range("F7").Value = 1551
$ws2.Range("F12").Value = 69

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 115

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 871
$ws4.Range("F4").Value = 13859
$ws4.Range("F5").Value = 13647
$ws4.Range("F6").Value = 1055
$ws4.Range("F9").Value = 605
$ws4.Range("F13").Value = 772
$ws4.Range("F16").Value = 2154
$ws4.Range("F17").Value = 117
$ws4.Range("F18").Value = 94
$ws4.Range("F19").Value = 80
$ws4.Range("F20").Value = 130
$ws4.Range("F24").Value = 115
$ws4.Range("F25").Value = 115
$ws4.Range("F26").Value = 538
$ws4.Range("F27").Value = 438
$ws4.Range("F28").Value = 417
$ws4.Range("F29").Value = 331
$ws4.Range("F31").Value = 845
$ws4.Range("F32").Value = 169
$ws4.Range("F33").Value = 1551
$ws4.Range("F38").Value = 100
$ws4.Range("F39").Value = 69
